$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.849.04"
$ws.Range("E2").Value = "  +1.04%  "

$ws.Range("D3").Value = "'1.767.16"
$ws.Range("E3").Value = "  +0.90%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  -0.38%  "

$ws.Range("D5").Value = "'327.44"
$ws.Range("E5").Value = "  +1.02%  "

$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "  -0.34%  "

$ws.Range("D7").Value = "'0.4484"
$ws.Range("E7").Value = "  -1.89%  "

$ws.Range("D8").Value = "'0.3548"
$ws.Range("E8").Value = "  -0.85%  "

$ws.Range("E9").Value = "  -0.68%  "

$ws.Range("D10").Value = "'42.11"
$ws.Range("E10").Value = "  +0.03%  "

$ws.Range("E11").Value = "  -0.10%  "

$ws.Range("E12").Value = "  -0.33%  "

$ws.Range("E13").Value = "  +0.56%  "

$ws.Range("D14").Value = "'6.020"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").Value = "'7.189"
$ws.Range("E15").Value = "  +1.55%  "

$ws.Range("D16").Value = "'1.765.00"
$ws.Range("E16").Value = "  +0.51%  "

$ws.Range("D17").Value = "'93.01"
$ws.Range("E17").Value = "  +0.76%  "

$ws.Range("D18").Value = "'0.00001055"
$ws.Range("E18").Value = "  -0.70%  "

$ws.Range("D19").Value = "'0.06431"
$ws.Range("E19").Value = "  +0.15%  "

$ws.Range("D20").Value = "'1.000"
$ws.Range("E20").Value = "  -0.27%  "

$ws.Range("E21").Value = "  +2.28%  "

$ws.Range("E22").Value = "  -0.74%  "

$ws.Range("D23").Value = "'27.890.19"
$ws.Range("E23").Value = "  +0.92%  "

$ws.Range("D24").Value = "'11.28"
$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").Value = "'2.111"
$ws.Range("E25").Value = "  -0.33%  "

$ws.Range("D26").Value = "'162.45"
$ws.Range("E26").Value = "  -1.03%  "

$ws.Range("D27").Value = "'20.20"
$ws.Range("E27").Value = "  -0.50%  "

$ws.Range("D28").Value = "'1.968.24"
$ws.Range("E28").Value = "  +0.60%  "

$ws.Range("E29").Value = "  +3.99%  "

$ws.Range("D30").Value = "'125.07"
$ws.Range("E30").Value = "  -0.94%  "

$ws.Range("D31").Value = "'1.091"

$ws.Range("D32").Value = "'0.09135"
$ws.Range("E32").Value = "  -0.47%  "

$ws.Range("D33").Value = "'3.655"
$ws.Range("E33").Value = "  -0.44%  "

$ws.Range("D34").Value = "'5.561"
$ws.Range("E34").Value = "  +0.59%  "

$ws.Range("E35").Value = "  -0.20%  "

$ws.Range("E36").Value = "  -0.39%  "

$ws.Range("D37").Value = "'0.06098"
$ws.Range("E37").Value = "  +0.90%  "

$ws.Range("D38").Value = "'0.2092"
$ws.Range("E38").Value = "  -0.28%  "

$ws.Range("D39").Value = "'4.958"
$ws.Range("E39").Value = "  -0.18%  "

$ws.Range("E40").Value = "  -0.83%  "

$ws.Range("E41").Value = "  -2.17%  "

$ws.Range("E42").Value = "  +0.68%  "

$ws.Range("D43").Value = "'7.917"
$ws.Range("E43").Value = "  +1.69%  "

$ws.Range("D44").Value = "'13.22"
$ws.Range("E44").Value = "  -0.35%  "

$ws.Range("D45").Value = "'3.738"
$ws.Range("E45").Value = "  +0.62%  "

$ws.Range("D46").Value = "'0.5847"
$ws.Range("E46").Value = "  -0.91%  "

$ws.Range("D47").Value = "'122.30"
$ws.Range("E47").Value = "  -0.73%  "

$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("D49").Value = "'0.06908"
$ws.Range("E49").Value = "  +0.67%  "

$ws.Range("D50").Value = "'1.134"
$ws.Range("E50").Value = "  -0.82%  "

$ws.Range("D51").Value = "'72.66"
$ws.Range("E51").Value = "  +0.94%  "
